$d = $word.ActiveDocument

# Position a zero-length range at the very end of the document body
# (after the last existing paragraph, before the section break) so
# that InsertXML appends rather than overwrites existing content.
$endPos = $d.Content.End
$r = $d.Range($endPos, $endPos)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# New empty paragraph (Normal style, only language mark carried over)
$emptyPara = '<w:p ' + $ns + '><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr></w:p>'

# New paragraph with the run.py execution path note
$textPara = '<w:p ' + $ns + '><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
            '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Ruta de Ejecución desde run.py</w:t></w:r></w:p>'

$null = $r.InsertXML($emptyPara + $textPara)
